# week 4 after-match inputs
# Adds a new "Week 39" column (AN) with after-match counts for the players
# who played that week, and clears a stray/unused border style that had
# been left on AK7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("AN1").Value = "Week 39"

# New week's data for the players who reported
$ws.Range("AN2").Value = 2
$ws.Range("AN4").Value = 7
$ws.Range("AN6").Value = 10
$ws.Range("AN8").Value = 3.5
$ws.Range("AN9").Value = 1.5

# AK7 had a stray border/fill style applied with no visible effect;
# clear it back to the default (no border) formatting.
$ws.Range("AK7").Borders.LineStyle = -4142

# Restore the frozen-pane column split and move the active selection to
# reflect where the user was last working.
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws.Range("B1").Select()
$aw.FreezePanes = $true
$ws.Range("AL6").Select()
